# Fruta / hortaliza, semanal
# Insert a new weekly data row above the current row 6 (existing rows 6-41
# shift down to 7-42) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push the existing data rows (6:41) down by one row.
$ws.Rows.Item(6).EntireRow.Insert()

# Populate the newly inserted row 6 with the new weekly record.
$ws.Cells.Item(6, 1).Value = 8
$ws.Cells.Item(6, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(6, 3).Value = "Coquimbo"
$ws.Cells.Item(6, 4).Value = 44881
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = 100112026
$ws.Cells.Item(6, 7).Value = "Haba"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 500
$ws.Cells.Item(6, 11).Value = 7000
$ws.Cells.Item(6, 12).Value = 8000
$ws.Cells.Item(6, 13).Value = 7500
$ws.Cells.Item(6, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(6, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(6, 16).Value = 300
$ws.Cells.Item(6, 17).Value = 25
$ws.Cells.Item(6, 18).Value = "Hortaliza"
